$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = '@'
    $rng.Value = $val
    $rng.Style = 'Normal'
}

Set-TextCell 'D2' '68.584.07'
Set-TextCell 'E2' '  -0.77%  '
Set-TextCell 'D3' '3.909.24'
Set-TextCell 'E3' '  +2.76%  '
Set-TextCell 'D4' '0.999'
Set-TextCell 'E4' '  -0.21%  '
Set-TextCell 'D5' '602.80'
Set-TextCell 'E5' '  +0.12%  '
Set-TextCell 'D6' '165.38'
Set-TextCell 'E6' '  +1.27%  '
Set-TextCell 'D7' '3.907.20'
Set-TextCell 'E7' '  +2.74%  '
Set-TextCell 'E8' '  -0.08%  '
Set-TextCell 'D9' '0.529'
Set-TextCell 'E9' '  -1.43%  '
Set-TextCell 'E10' '  -2.99%  '
Set-TextCell 'D11' '6.38'
Set-TextCell 'E11' '  +1.23%  '
Set-TextCell 'E12' '  +0.17%  '
Set-TextCell 'D13' '37.28'
Set-TextCell 'E13' '  -0.03%  '
Set-TextCell 'D14' '0.0000248'
Set-TextCell 'E14' '  +0.69%  '
Set-TextCell 'D15' '4.564.14'
Set-TextCell 'E15' '  +2.77%  '
Set-TextCell 'D16' '3.923.50'
Set-TextCell 'E16' '  +3.23%  '
Set-TextCell 'D17' '68.647.49'
Set-TextCell 'E17' '  -0.84%  '
Set-TextCell 'E18' '  +0.11%  '
Set-TextCell 'D19' '17.21'
Set-TextCell 'E19' '  -0.77%  '
Set-TextCell 'E20' '  -1.33%  '
Set-TextCell 'D21' '11.03'
Set-TextCell 'E21' '  -2.11%  '
Set-TextCell 'D22' '486.19'
Set-TextCell 'E22' '  -0.88%  '
Set-TextCell 'D23' '0.725'
Set-TextCell 'E23' '  +0.36%  '
Set-TextCell 'D24' '0.0000168'
Set-TextCell 'E24' '  +10.72%  '
Set-TextCell 'D25' '84.70'
Set-TextCell 'E25' '  +0.01%  '
Set-TextCell 'E26' '  -0.99%  '
Set-TextCell 'D27' '12.08'
Set-TextCell 'E27' '  -1.24%  '
Set-TextCell 'D28' '10.12'
Set-TextCell 'E28' '  +0.60%  '
Set-TextCell 'E29' '  -0.05%  '
Set-TextCell 'E30' '  -1.62%  '
Set-TextCell 'D31' '4.061.26'
Set-TextCell 'E31' '  +2.76%  '
Set-TextCell 'E32' '  -0.67%  '
Set-TextCell 'E33' '  -3.67%  '
Set-TextCell 'D34' '31.97'
Set-TextCell 'E34' '  -0.10%  '
Set-TextCell 'D35' '3.857.16'
Set-TextCell 'E35' '  +2.82%  '
Set-TextCell 'E36' '  +0.00%  '
Set-TextCell 'E37' '  +2.18%  '
Set-TextCell 'B38' 'Filecoin'
Set-TextCell 'C38' 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextCell 'D38' '5.91'
Set-TextCell 'E38' '  +0.04%  '
Set-TextCell 'B39' 'Kaspa'
Set-TextCell 'C39' 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextCell 'D39' '0.138'
Set-TextCell 'E39' '  -1.83%  '
Set-TextCell 'D40' '0.999'
Set-TextCell 'E40' '  -0.17%  '
Set-TextCell 'D41' '3.14'
Set-TextCell 'E41' '  +3.24%  '
Set-TextCell 'D42' '0.317'
Set-TextCell 'E42' '  -2.06%  '
Set-TextCell 'D43' '427.43'
Set-TextCell 'E43' '  +1.29%  '
Set-TextCell 'D44' '48.28'
Set-TextCell 'E44' '  -0.43%  '
Set-TextCell 'E45' '  -0.38%  '
Set-TextCell 'E46' '  +1.21%  '
Set-TextCell 'E47' '  -0.01%  '
Set-TextCell 'D48' '141.92'
Set-TextCell 'E48' '  +0.03%  '
Set-TextCell 'D49' '2.815.55'
Set-TextCell 'E49' '  -0.42%  '
Set-TextCell 'D50' '26.09'
Set-TextCell 'E50' '  +5.43%  '
Set-TextCell 'B51' 'Arweave'
Set-TextCell 'C51' 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
Set-TextCell 'D51' '39.15'
Set-TextCell 'E51' '  -1.00%  '
